$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the free-text summary cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.79 = 10434.48 pesos`n✅ 10434.48 pesos = 2.77 = 943.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 357.99
$ws2.Range("O10").Value = 3735.44
$ws2.Range("N12").Value = 3761.9
$ws2.Range("O12").Value = 340
